$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2511715
$ws.Range("I43").Value = 4053274.5
$ws.Range("K43").Value = 4053274.5
$ws.Range("M43").Value = -4053205.5

$ws.Range("H96").Value = 1747.6666
$ws.Range("I96").Value = 1877.2
$ws.Range("K96").Value = 5631.6
$ws.Range("M96").Value = -4258.6

$ws.Range("H113").Value = 7374.3335
$ws.Range("I113").Value = 2956.1428
$ws.Range("K113").Value = 2956.1428
$ws.Range("M113").Value = 297.8571999999999

$ws.Range("H137").Value = 4830
$ws.Range("J137").Value = 4000
$ws.Range("L137").Value = 12000
$ws.Range("N137").Value = -17100

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 3902.7083
$ws.Range("I97").Value = 1338.875
$ws.Range("K97").Value = 1338.875
$ws.Range("M97").Value = -842.875

$ws.Range("H102").Value = 1718.8182
$ws.Range("I102").Value = 755.1111
$ws.Range("J102").Value = 6055.5
$ws.Range("K102").Value = 755.1111
$ws.Range("L102").Value = 6055.5
$ws.Range("M102").Value = 866.8889
$ws.Range("N102").Value = -9299.5

$ws.Range("H110").Value = 3679.6667
$ws.Range("I110").Value = 2684.1667
$ws.Range("K110").Value = 2684.1667
$ws.Range("M110").Value = -639.1667000000002

$ws.Range("H122").Value = 3817.5
$ws.Range("I122").Value = 4232
$ws.Range("J122").Value = 2988.5
$ws.Range("K122").Value = 12696
$ws.Range("L122").Value = 8965.5
$ws.Range("M122").Value = -10246
$ws.Range("N122").Value = -13865.5

$ws.Range("H132").Value = 7034.514
$ws.Range("I132").Value = 7263.909
$ws.Range("K132").Value = 21791.727
$ws.Range("M132").Value = -19261.727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 887.625
$ws.Range("I22").Value = 887.625
$ws.Range("K22").Value = 887.625
$ws.Range("M22").Value = -714.625

$ws.Range("H26").Value = 15161.667
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

$ws.Range("H94").Value = 437.6154
$ws.Range("I94").Value = 460.36365
$ws.Range("K94").Value = 460.36365
$ws.Range("M94").Value = -9.363650000000007

$ws.Range("H99").Value = 5307.7144
$ws.Range("J99").Value = 6383.3335
$ws.Range("L99").Value = 6383.3335
$ws.Range("N99").Value = -9379.333500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2856.4546
$ws.Range("I16").Value = 2692.1
$ws.Range("K16").Value = 2692.1
$ws.Range("M16").Value = -2405.1

$ws.Range("H105").Value = 5428
$ws.Range("I105").Value = 5231.5
$ws.Range("J105").Value = 7000
$ws.Range("K105").Value = 5231.5
$ws.Range("L105").Value = 7000
$ws.Range("M105").Value = -3484.5
$ws.Range("N105").Value = -10494

$ws.Range("H113").Value = 2856.4546
$ws.Range("I113").Value = 2692.1
$ws.Range("K113").Value = 2692.1
$ws.Range("M113").Value = -522.0999999999999

$ws.Range("H122").Value = 6452.4614
$ws.Range("I122").Value = 5442.5454
$ws.Range("J122").Value = 12007
$ws.Range("K122").Value = 16327.6362
$ws.Range("L122").Value = 36021
$ws.Range("M122").Value = -13877.6362
$ws.Range("N122").Value = -40921

$ws.Range("H132").Value = 1373
$ws.Range("I132").Value = 955
$ws.Range("K132").Value = 2865
$ws.Range("M132").Value = -335

$ws.Range("H134").Value = 1199.6666
$ws.Range("I134").Value = 1199.6666
$ws.Range("K134").Value = 3598.9998
$ws.Range("M134").Value = -1063.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 168.52
$ws.Range("I2").Value = 129
$ws.Range("K2").Value = 774
$ws.Range("M2").Value = -661

$ws.Range("H12").Value = 644.4375
$ws.Range("I12").Value = 461.6
$ws.Range("J12").Value = 949.1667
$ws.Range("K12").Value = 1384.8
$ws.Range("L12").Value = 2847.5001
$ws.Range("M12").Value = -1211.8
$ws.Range("N12").Value = -3193.5001

$ws.Range("H23").Value = 904.25
$ws.Range("I23").Value = 16.5
$ws.Range("J23").Value = 1031.0714
$ws.Range("K23").Value = 49.5
$ws.Range("L23").Value = 3093.2142
$ws.Range("M23").Value = 185.5
$ws.Range("N23").Value = -3563.2142

$ws.Range("H138").Value = 1007058.3
$ws.Range("I138").Value = 1675097.1
$ws.Range("K138").Value = 5025291.300000001
$ws.Range("M138").Value = -5020151.300000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 170.625
$ws.Range("I2").Value = 223.83333
$ws.Range("J2").Value = 11
$ws.Range("K2").Value = 223.83333
$ws.Range("L2").Value = 11
$ws.Range("M2").Value = -110.83333
$ws.Range("N2").Value = -237

$ws.Range("H42").Value = 80000
$ws.Range("J42").Value = 80000
$ws.Range("L42").Value = 80000
$ws.Range("N42").Value = -80970

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H70").Value = 13491.786
$ws.Range("I70").Value = 6747.8335
$ws.Range("J70").Value = 18549.75
$ws.Range("K70").Value = 6747.8335
$ws.Range("L70").Value = 18549.75
$ws.Range("M70").Value = -6477.8335
$ws.Range("N70").Value = -19089.75

$ws.Range("H73").Value = 13491.786
$ws.Range("I73").Value = 6747.8335
$ws.Range("J73").Value = 18549.75
$ws.Range("K73").Value = 6747.8335
$ws.Range("L73").Value = 18549.75
$ws.Range("M73").Value = -5811.8335
$ws.Range("N73").Value = -20421.75

$ws.Range("H113").Value = 2620.9167
$ws.Range("I113").Value = 2151.6
$ws.Range("K113").Value = 2151.6
$ws.Range("M113").Value = 18.40000000000009

$ws.Range("H115").Value = 80000
$ws.Range("J115").Value = 80000
$ws.Range("L115").Value = 80000
$ws.Range("N115").Value = -82350

$ws.Range("H122").Value = 125004296
$ws.Range("I122").Value = 166670450
$ws.Range("J122").Value = 5847.5
$ws.Range("K122").Value = 500011350
$ws.Range("L122").Value = 17542.5
$ws.Range("M122").Value = -500008900
$ws.Range("N122").Value = -22442.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9044.6
$ws.Range("I40").Value = 9044.6
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 9044.6
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -8908.6
$ws.Range("N40").ClearContents()

$ws.Range("H82").Value = 7367.9414
$ws.Range("I82").Value = 17729.334
$ws.Range("K82").Value = 17729.334
$ws.Range("M82").Value = -17368.334

$ws.Range("H85").Value = 7367.9414
$ws.Range("I85").Value = 17729.334
$ws.Range("K85").Value = 17729.334
$ws.Range("M85").Value = -16481.334

$ws.Range("H93").Value = 3238.5
$ws.Range("I93").Value = 3001
$ws.Range("J93").Value = 3476
$ws.Range("K93").Value = 3001
$ws.Range("L93").Value = 3476
$ws.Range("M93").Value = -1753
$ws.Range("N93").Value = -5972

$ws.Range("H136").Value = 4711.564
$ws.Range("I136").Value = 4793.3716
$ws.Range("K136").Value = 14380.1148
$ws.Range("M136").Value = -11830.1148

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 6310.304
$ws.Range("I4").Value = 7194.684
$ws.Range("K4").Value = 7194.684
$ws.Range("M4").Value = -7081.684

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0

$ws.Range("H81").Value = 5585.5713
$ws.Range("I81").Value = 5585.5713
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 11171.1426
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -10110.1426
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 5585.5713
$ws.Range("I84").Value = 5585.5713
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 55855.713
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -50551.713
$ws.Range("N84").ClearContents()

$ws.Range("H122").Value = 4479.4375
$ws.Range("I122").Value = 5205.769
$ws.Range("J122").Value = 1332
$ws.Range("K122").Value = 15617.307
$ws.Range("L122").Value = 4196
$ws.Range("M122").Value = -13167.307
$ws.Range("N122").Value = -8896
